$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 13 (sample with no fastq file) - shifts rows 14:37 up to 13:36
$ws.Rows("13:13").Delete()

# Update the active selection to match the author's post-edit view
# (the row that used to be row 14 is now row 13, selected as a whole row)
$ws.Rows("13:13").Select()
